# Add a new bold paragraph after the "... múlt heti sorsolást”." sentence,
# describing the May 31 update, matching the author's commit.
#
# Strategy: locate the end of the existing last sentence with Find (so the
# script doesn't depend on hard-coded character offsets), then inject a
# brand-new <w:p> (paragraph break + 3 runs, with proofErr spell markers
# around "ben") right at that point via Range.InsertXML, using real OOXML
# markup equivalent to what Word itself would produce.

$d = $word.ActiveDocument

$find = $d.Content.Find
$found = $find.Execute("múlt heti sorsolást”.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor sentence not found in document."
}

$insertAt = $find.Parent.End
$ins = $d.Range($insertAt, $insertAt)

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Május 31-én szintén nem volt sok haladás, továbbra is az Excel-</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>ben</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> lévő adatok kiírása történt, csak abba volt haladás az órán.</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$ins.InsertXML($newParaXml)

Write-Host "Inserted new paragraph after anchor; paragraph count now:" $d.Paragraphs.Count
